# Auto-generated edit script applying the diff to Behemoth_Profits workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 1760
$ws.Range("I41").Value = 1620.3
$ws.Range("K41").Value = 1620.3
$ws.Range("M41").Value = -1180.3

$ws.Range("H62").Value = 16762.5
$ws.Range("I62").Value = 2050
$ws.Range("K62").Value = 2050
$ws.Range("M62").Value = -1426

$ws.Range("H64").Value = 4961.5
$ws.Range("I64").Value = 4000
$ws.Range("K64").Value = 4000
$ws.Range("M64").Value = -3752

$ws.Range("H65").Value = 16762.5
$ws.Range("I65").Value = 2050
$ws.Range("K65").Value = 10250
$ws.Range("M65").Value = -7130

$ws.Range("H67").Value = 4961.5
$ws.Range("I67").Value = 4000
$ws.Range("K67").Value = 4000
$ws.Range("M67").Value = -3142

$ws.Range("H70").Value = 1539.7931
$ws.Range("J70").Value = 1429.6111
$ws.Range("L70").Value = 4288.8333
$ws.Range("N70").Value = -4828.8333

$ws.Range("H73").Value = 1539.7931
$ws.Range("J73").Value = 1429.6111
$ws.Range("L73").Value = 4288.8333
$ws.Range("N73").Value = -6160.8333

$ws.Range("H101").Value = 2356.2
$ws.Range("I101").Value = 2153
$ws.Range("J101").Value = 4185
$ws.Range("K101").Value = 6459
$ws.Range("L101").Value = 12555
$ws.Range("M101").Value = -4837
$ws.Range("N101").Value = -15799

$ws.Range("H113").Value = 38464188
$ws.Range("I113").Value = 14287764
$ws.Range("J113").Value = 66670016
$ws.Range("K113").Value = 14287764
$ws.Range("L113").Value = 66670016
$ws.Range("M113").Value = -14284510
$ws.Range("N113").Value = -66676524

$ws.Range("H127").Value = 2710
$ws.Range("I127").Value = 2710
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 8130
$ws.Range("L127").Value = 0
$ws.Range("M127").Value = -3170
$ws.Range("N127").ClearContents()

$ws.Range("H132").Value = 2931.8
$ws.Range("I132").Value = 3044
$ws.Range("J132").Value = 800
$ws.Range("K132").Value = 9132
$ws.Range("L132").Value = 2400
$ws.Range("M132").Value = -6602
$ws.Range("N132").Value = -7460

$ws.Range("H137").Value = 5042
$ws.Range("I137").Value = 1784
$ws.Range("J137").Value = 7834.5713
$ws.Range("K137").Value = 5352
$ws.Range("L137").Value = 23503.7139
$ws.Range("M137").Value = -2802
$ws.Range("N137").Value = -28603.7139

$ws.Range("H138").Value = 2341.157
$ws.Range("J138").Value = 3139.8667
$ws.Range("L138").Value = 9419.6001
$ws.Range("N138").Value = -19699.6001

$ws.Range("H141").Value = 1708.4762
$ws.Range("I141").Value = 1708.4762
$ws.Range("K141").Value = 5125.4286
$ws.Range("M141").Value = 54.57139999999981

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2803.7693
$ws.Range("I45").Value = 2444.9
$ws.Range("J45").Value = 4000
$ws.Range("K45").Value = 2444.9
$ws.Range("L45").Value = 4000
$ws.Range("M45").Value = -2067.9
$ws.Range("N45").Value = -4754

$ws.Range("H63").Value = 6179.2
$ws.Range("I63").Value = 4268.8
$ws.Range("J63").Value = 10000
$ws.Range("K63").Value = 4268.8
$ws.Range("L63").Value = 10000
$ws.Range("M63").Value = -3582.8
$ws.Range("N63").Value = -11372

$ws.Range("H66").Value = 6179.2
$ws.Range("I66").Value = 4268.8
$ws.Range("J66").Value = 10000
$ws.Range("K66").Value = 21344
$ws.Range("L66").Value = 50000
$ws.Range("M66").Value = -17912
$ws.Range("N66").Value = -56864

$ws.Range("H110").Value = 1580
$ws.Range("I110").Value = 1580
$ws.Range("K110").Value = 1580
$ws.Range("M110").Value = 465

$ws.Range("H122").Value = 2313.611
$ws.Range("I122").Value = 1220.4166
$ws.Range("J122").Value = 4500
$ws.Range("K122").Value = 3661.2498
$ws.Range("L122").Value = 13500
$ws.Range("M122").Value = -1211.2498
$ws.Range("N122").Value = -18400

$ws.Range("H132").Value = 5353.5186
$ws.Range("I132").Value = 2043.25
$ws.Range("K132").Value = 6129.75
$ws.Range("M132").Value = -3599.75

$ws.Range("H139").Value = 80580.836
$ws.Range("J139").Value = 80580.836
$ws.Range("L139").Value = 80580.836
$ws.Range("N139").Value = -90860.836

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 51851
$ws.Range("J81").Value = 51851
$ws.Range("L81").Value = 51851
$ws.Range("N81").Value = -53973

$ws.Range("H84").Value = 51851
$ws.Range("J84").Value = 51851
$ws.Range("L84").Value = 155553
$ws.Range("N84").Value = -166161

$ws.Range("H94").Value = 1324.5
$ws.Range("I94").Value = 1563.9166
$ws.Range("J94").Value = 845.6667
$ws.Range("K94").Value = 1563.9166
$ws.Range("L94").Value = 845.6667
$ws.Range("M94").Value = -1112.9166
$ws.Range("N94").Value = -1747.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 66923.336
$ws.Range("J20").Value = 66923.336
$ws.Range("L20").Value = 66923.336
$ws.Range("N20").Value = -67395.336

$ws.Range("H30").Value = 66923.336
$ws.Range("J30").Value = 66923.336
$ws.Range("L30").Value = 66923.336
$ws.Range("N30").Value = -67105.336

$ws.Range("H31").Value = 605219.3
$ws.Range("I31").Value = 6823.6665
$ws.Range("J31").Value = 1303347.5
$ws.Range("K31").Value = 6823.6665
$ws.Range("L31").Value = 1303347.5
$ws.Range("M31").Value = -6528.6665
$ws.Range("N31").Value = -1303937.5

$ws.Range("H34").Value = 605219.3
$ws.Range("I34").Value = 6823.6665
$ws.Range("J34").Value = 1303347.5
$ws.Range("K34").Value = 6823.6665
$ws.Range("L34").Value = 1303347.5
$ws.Range("M34").Value = -6621.6665
$ws.Range("N34").Value = -1303751.5

$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

$ws.Range("H99").Value = 4832.75
$ws.Range("I99").Value = 4837.4287
$ws.Range("K99").Value = 4837.4287
$ws.Range("M99").Value = -3339.4287

$ws.Range("H126").Value = 4832.75
$ws.Range("I126").Value = 4837.4287
$ws.Range("K126").Value = 14512.2861
$ws.Range("M126").Value = -12042.2861

$ws.Range("H128").Value = 66923.336
$ws.Range("J128").Value = 66923.336
$ws.Range("L128").Value = 66923.336
$ws.Range("N128").Value = -76883.336

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()

$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()

$ws.Range("H125").Value = 92999.664
$ws.Range("J125").Value = 92999.664
$ws.Range("L125").Value = 92999.664
$ws.Range("N125").Value = -97919.664

$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()

$ws.Range("H132").Value = 62503372
$ws.Range("I132").Value = 90912580
$ws.Range("J132").Value = 3110.2
$ws.Range("K132").Value = 272737740
$ws.Range("L132").Value = 9330.599999999999
$ws.Range("M132").Value = -272735210
$ws.Range("N132").Value = -14390.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H99").Value = 58366.832
$ws.Range("I99").Value = 51552.75
$ws.Range("J99").Value = 71995
$ws.Range("K99").Value = 51552.75
$ws.Range("L99").Value = 71995
$ws.Range("M99").Value = -48557.75
$ws.Range("N99").Value = -77985

$ws.Range("H132").Value = 56510.684
$ws.Range("I132").Value = 4553.857
$ws.Range("K132").Value = 13661.571
$ws.Range("M132").Value = -11131.571

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H37").Value = 30495
$ws.Range("I37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("M37").ClearContents()

$ws.Range("H100").Value = 1512.1428
$ws.Range("I100").Value = 1512.1428
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 3024.2856
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -2483.2856
$ws.Range("N100").ClearContents()

$ws.Range("H122").Value = 13715.434
$ws.Range("I122").Value = 9474.75
$ws.Range("J122").Value = 22196.8
$ws.Range("K122").Value = 28424.25
$ws.Range("L122").Value = 66590.39999999999
$ws.Range("M122").Value = -25974.25
$ws.Range("N122").Value = -71490.39999999999

$ws.Range("H126").Value = 6839.174
$ws.Range("I126").Value = 6365.45
$ws.Range("J126").Value = 9997.333000000001
$ws.Range("K126").Value = 19096.35
$ws.Range("L126").Value = 29991.999
$ws.Range("M126").Value = -16626.35
$ws.Range("N126").Value = -34931.999

$ws.Range("H127").Value = 64000
$ws.Range("J127").Value = 64000
$ws.Range("L127").Value = 64000
$ws.Range("N127").Value = -73920

$ws.Range("H132").Value = 1887.28
$ws.Range("I132").Value = 1705.1177
$ws.Range("J132").Value = 2274.375
$ws.Range("K132").Value = 5115.3531
$ws.Range("L132").Value = 6823.125
$ws.Range("M132").Value = -2585.3531
$ws.Range("N132").Value = -11883.125
